$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" number (e.g. "581.78") must be forced
# back to Text so they match the sheets existing text-formatted Price column
# (values like "70.847.25" already fail numeric parsing and stay text on their own).

$ws.Range("D2").Value = "70.847.25"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "3.651.63"
$ws.Range("E3").Value = "  +6.66%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.78"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.18"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "3.639.40"
$ws.Range("E7").Value = "  +6.53%  "

$ws.Range("E8").Value = "  +3.59%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.198"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.86"
$ws.Range("E11").Value = "  +25.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.606"
$ws.Range("E12").Value = "  +3.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.63"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000290"
$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").Value = "4.240.57"
$ws.Range("E15").Value = "  +6.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "672.18"
$ws.Range("E16").Value = "  -2.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.96"

$ws.Range("D18").Value = "3.680.26"
$ws.Range("E18").Value = "  +7.49%  "

$ws.Range("D19").Value = "70.912.67"
$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.80"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.48"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.933"
$ws.Range("E23").Value = "  +4.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.16"
$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.12"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +5.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.87"
$ws.Range("E30").Value = "  +4.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.38"
$ws.Range("E31").Value = "  +4.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.03"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  -3.13%  "

$ws.Range("E34").Value = "  +5.79%  "

$ws.Range("E35").Value = "  +7.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "579.75"
$ws.Range("E36").Value = "  +1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.10"
$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("E38").Value = "  +4.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.61"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.596.28"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0456"
$ws.Range("E42").Value = "  +8.87%  "

$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.347"
$ws.Range("E44").Value = "  +4.63%  "

$ws.Range("D45").Value = "0.0₃0750"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "34.93"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"

$ws.Range("E49").Value = "  +3.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.82"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("E51").Value = "  +8.46%  "
